# Squash merge framework/feature-SNAPSHOT.20241005 into dev
#
# Adds two new localization resource rows (19 & 20) to the
# ResXResourceManager sheet for the new "ArgumentException_empty_string"
# and "ArgumentException_empty_or_whitespace" resource keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$project = "Disco.Localization.Resources"
$file    = "Strings"

# --- Row 19: ArgumentException_empty_string -------------------------------
$key19 = "ArgumentException_empty_string"
$en19  = 'The string parameter "{0}" is not allowed to be NULL or empty.'
$zh19  = '字符串参数“{0}”不允许为 NULL 或空白。'

$ws.Cells.Item(19, 1).Value = $project
$ws.Cells.Item(19, 2).Value = $file
$ws.Cells.Item(19, 3).Value = $key19
$ws.Cells.Item(19, 5).Value = $en19
$ws.Cells.Item(19, 7).Value = $en19
$ws.Cells.Item(19, 9).Value = $zh19

# --- Row 20: ArgumentException_empty_or_whitespace -------------------------
$key20 = "ArgumentException_empty_or_whitespace"
$en20  = 'The string parameter "{0}" is not allowed to be NULL or empty or white-spaces.'
$zh20  = '字符串参数“{0}”不允许为 NULL 或空格符。'

$ws.Cells.Item(20, 1).Value = $project
$ws.Cells.Item(20, 2).Value = $file
$ws.Cells.Item(20, 3).Value = $key20
$ws.Cells.Item(20, 5).Value = $en20
$ws.Cells.Item(20, 7).Value = $en20
$ws.Cells.Item(20, 9).Value = $zh20

# --- Match formatting used by the rest of the table (copy cell formats
# from the row above, which already carries the correct styles) -----------
$colsWithData = @(1, 2, 3, 5, 7, 9)
foreach ($col in $colsWithData) {
    $ws.Cells.Item(18, $col).Copy()
    $ws.Cells.Item(19, $col).PasteSpecial(-4122)
    $ws.Cells.Item(18, $col).Copy()
    $ws.Cells.Item(20, $col).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# The new rows use the same (non-default-measured) custom row height as the
# rest of the data rows in the table.
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(18).RowHeight

# --- Match the active selection recorded in the sheet (bottom-right pane
# ends up with I20 selected after adding the new rows) ----------------------
$ws.Range("I20").Select()
